# reporte_ventas.xlsx - add "ID Venta" and "Cantidad de Producto" columns,
# refresh the sales rows (consolidate the three "Juan Morales Morales"
# entries into a single, more recent one) and move the Total sum formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the two rows that disappear (three "Juan Morales Morales" rows
#    collapse into one), working bottom-up so row numbers stay valid.
#    This also pulls the old totals row (13, =SUM(E5:E12)) up to row 11
#    and keeps it correctly ranged to the surviving rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# 2) Headers (row 5) -- A:D/new-E already exist with the header style
#    (s="2"); F5/G5 are brand new cells, so copy the header style from
#    an existing header cell before writing into them.
# ---------------------------------------------------------------------
$ws.Range("D5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A5").Value = "ID Venta"
$ws.Range("B5").Value = "Fecha de Compra"
$ws.Range("C5").Value = "Nombres"
$ws.Range("D5").Value = "Tipo de Comprobante"
$ws.Range("E5").Value = "Cantidad de Producto"
$ws.Range("F5").Value = "Forma de Pago"
$ws.Range("G5").Value = "Total"

# ---------------------------------------------------------------------
# 3) Data rows 6-10 -- copy the data-row style (s="3") into the new F/G
#    cells before writing data so every new cell picks up the bordered
#    style used by the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("D6:D10").Copy()
$ws.Range("F6:F10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("E6:E10").Copy()
$ws.Range("G6:G10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 6
$ws.Range("A6").Value = 1001
$ws.Range("B6").Value = "2021-05-06 12:45:51"
$ws.Range("C6").Value = "Leonel Pérez Pérez"
$ws.Range("D6").Value = "Boleta"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "Efectivo"
$ws.Range("G6").Value = 5990

# Row 7
$ws.Range("A7").Value = 1002
$ws.Range("B7").Value = "2021-05-07 13:02:20"
$ws.Range("C7").Value = "Sergio Diaz Figueroa"
$ws.Range("D7").Value = "Factura"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "Debito"
$ws.Range("G7").Value = 221030

# Row 8
$ws.Range("A8").Value = 1003
$ws.Range("B8").Value = "2021-05-08 15:43:06"
$ws.Range("C8").Value = "Leonel Pérez Pérez"
$ws.Range("D8").Value = "Boleta"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "Debito"
$ws.Range("G8").Value = 10390

# Row 9
$ws.Range("A9").Value = 1004
$ws.Range("B9").Value = "2021-05-08 15:55:05"
$ws.Range("C9").Value = "Leonel Pérez Pérez"
$ws.Range("D9").Value = "Boleta"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = "Debito"
$ws.Range("G9").Value = 10790

# Row 10 (was 3 separate "Juan Morales Morales" rows, now consolidated)
$ws.Range("A10").Value = 1007
$ws.Range("B10").Value = "2021-06-06 22:55:24"
$ws.Range("C10").Value = "Juan Morales Morales"
$ws.Range("D10").Value = "Boleta"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "Debito"
$ws.Range("G10").Value = 27890

# ---------------------------------------------------------------------
# 4) Totals formula now lives in G11 (row 11 already holds the shifted
#    =SUM(E5:E10) from step 1 in column E -- move it over to G and widen
#    the summed range to match the new table).
# ---------------------------------------------------------------------
$ws.Range("E11").ClearContents()
$ws.Range("G11").Formula = "=SUM(G5:G10)"

# ---------------------------------------------------------------------
# 5) Column widths (A:G) per the new layout. ColumnWidth is stored with
#    a constant +5/6 character offset by this engine, so back it out to
#    land on the exact target width.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = (12 - 5/6)
$ws.Columns.Item(2).ColumnWidth = (20 - 5/6)
$ws.Columns.Item(3).ColumnWidth = (30 - 5/6)
$ws.Columns.Item(4).ColumnWidth = (25 - 5/6)
$ws.Columns.Item(5).ColumnWidth = (25 - 5/6)
$ws.Columns.Item(6).ColumnWidth = (20 - 5/6)
$ws.Columns.Item(7).ColumnWidth = (10 - 5/6)

# ---------------------------------------------------------------------
# 6) Selection / active cell to match the new used range.
# ---------------------------------------------------------------------
$ws.Range("A5:G10").Select()
$ws.Range("A5").Activate()
